$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")
$ws.Activate()

# New rows 37-39 describe additional key-condition lookups (acdtl*).
# Column B uses the wrap-text style already used by the rest of column B
# in this table (style index 3 in the original file).
$ws.Range("B37").WrapText = $true
$ws.Range("B38").WrapText = $true
$ws.Range("B39").WrapText = $true

# Row 37
$ws.Range("B37").Value = "AcDate = ,AND RelTxseq = "
$ws.Range("A37").Value = "acdtlRelTxseqEq2"

# Row 38
$ws.Range("A38").Value = "acdtlSlipNo2"
$ws.Range("B38").Value = "AcDate = ,AND SlipNo = "

# Row 39
$ws.Range("B39").Value = "RelDy = ,AND SlipNo = "
$ws.Range("A39").Value = "acdtlSlipNo"

# Column C repeats the existing sort-order note for each new row.
$ws.Range("C37").Value = "AcSeq ASC"
$ws.Range("C38").Value = "AcSeq ASC"
$ws.Range("C39").Value = "AcSeq ASC"

# Move the selection to the cell below the newly added rows, matching
# where the user ended up after typing the new entries.
$ws.Range("A40").Select() | Out-Null
